$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'75.983.17"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "'2.892.17"
$ws.Range("E3").Value = "  +6.49%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'195.51"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").Value = "'598.03"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D9").Value = "'0.193"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "'2.889.95"
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("D11").Value = "'0.399"
$ws.Range("E11").Value = "  +10.29%  "
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").Value = "'3.415.71"
$ws.Range("E14").Value = "  +6.23%  "
$ws.Range("D15").Value = "'75.847.68"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'0.0000191"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "'27.39"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "'2.874.05"
$ws.Range("E18").Value = "  +5.42%  "
$ws.Range("D19").Value = "'8.91"
$ws.Range("E19").Value = "  -4.52%  "
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("D21").Value = "'377.84"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'4.16"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("D24").Value = "'71.48"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "'3.027.05"
$ws.Range("E26").Value = "  +5.95%  "
$ws.Range("D27").Value = "'4.24"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").Value = "'9.84"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "'0.0000108"
$ws.Range("E29").Value = "  +9.21%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'509.37"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").Value = "'7.80"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'20.24"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "'163.28"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("D40").Value = "'182.91"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'0.346"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "'5.01"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("E45").Value = "  +7.27%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'40.37"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  -0.96%  "
$ws.Range("D49").Value = "'0.580"
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("D50").Value = "'3.77"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "'0.667"
$ws.Range("E51").Value = "  +11.67%  "
